$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as text so values like
# "0.7410" or "29.936.51" are preserved exactly instead of being
# auto-converted to numbers (which would drop formatting/trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.936.51"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.33"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7410"
$ws.Range("E5").Value = "  -3.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.46"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3155"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07189"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.75"
$ws.Range("E10").Value = "  -3.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08472"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7531"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.903.01"
$ws.Range("E13").Value = "  -8.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.397"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.69"
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.130"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.955.49"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.61"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.34"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007825"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.123.19"
$ws.Range("E22").Value = "  -10.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.986"
$ws.Range("E23").Value = "  -2.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1560"
$ws.Range("E25").Value = "  -2.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.312"
$ws.Range("E26").Value = "  -2.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.35"
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.64"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.479"
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.607"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.528"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.292"
$ws.Range("E33").Value = "  +4.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05332"
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.240"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7561"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9970"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.694"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01959"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.751"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4482"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.113.28"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.111"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.57"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8596"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.14"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.678"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.843"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.069"
$ws.Range("E50").Value = "  +3.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.022.29"
$ws.Range("E51").Value = "  -7.52%  "
